# Shrink the font on two flow-diagram boxes on slide 1:
#   "ADC/FFT"  (Rounded Rectangle 3)  : 8pt   (sz="800") -> 7pt    (sz="700")
#   "SIG"      (Rounded Rectangle 13) : 9.35pt(sz="935") -> 7.35pt (sz="735")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }

    $tr = $shp.TextFrame.TextRange
    $text = $tr.Text

    if ($text -eq "ADC/FFT") {
        $tr.Font.Size = 7
    } elseif ($text -eq "SIG") {
        $tr.Font.Size = 7.35
    }
}
